# Auto-generated script to apply 2025-05-15 data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2272
$ws.Range('L3').Value = 2295
$ws.Range('L4').Value = 625
$ws.Range('I5').Value = 728
$ws.Range('L5').Value = 137
$ws.Range('L6').Value = 2069
$ws.Range('I7').Value = 26303
$ws.Range('L7').Value = 7398

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L6').Value = 30
$ws.Range('L7').Value = 93

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L3').Value = 159
$ws.Range('L4').Value = 35
$ws.Range('L5').Value = 20
$ws.Range('L6').Value = 122
$ws.Range('L7').Value = 470

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L3').Value = 73
$ws.Range('L5').Value = 4
$ws.Range('L7').Value = 181

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 88
$ws.Range('L3').Value = 107
$ws.Range('L7').Value = 335

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L3').Value = 80
$ws.Range('L6').Value = 82
$ws.Range('L7').Value = 271

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L4').Value = 6
$ws.Range('L7').Value = 141

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L6').Value = 25
$ws.Range('L7').Value = 117

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L3').Value = 6
$ws.Range('L7').Value = 239
$ws.Range('L8').Value = 470
$ws.Range('L10').Value = 48
$ws.Range('L11').Value = 131
$ws.Range('L13').Value = 9
$ws.Range('L14').Value = 39
$ws.Range('L15').Value = 51
$ws.Range('L18').Value = 53
$ws.Range('L20').Value = 188
$ws.Range('L23').Value = 75
$ws.Range('L25').Value = 41
$ws.Range('L27').Value = 76
$ws.Range('L29').Value = 376
$ws.Range('L31').Value = 72
$ws.Range('L33').Value = 335
$ws.Range('L37').Value = 271
$ws.Range('L40').Value = 17
$ws.Range('L42').Value = 232
$ws.Range('L43').Value = 59
$ws.Range('L44').Value = 53
$ws.Range('L48').Value = 101
$ws.Range('L49').Value = 42
$ws.Range('L51').Value = 84
$ws.Range('L53').Value = 93
$ws.Range('L54').Value = 152
$ws.Range('L55').Value = 66
$ws.Range('L57').Value = 32
$ws.Range('I63').Value = 258
$ws.Range('L63').Value = 21
$ws.Range('L65').Value = 141
$ws.Range('L67').Value = 265
$ws.Range('L73').Value = 58
$ws.Range('L78').Value = 101
$ws.Range('L79').Value = 204
$ws.Range('L83').Value = 181
$ws.Range('L84').Value = 73
$ws.Range('L85').Value = 386
$ws.Range('L86').Value = 57
$ws.Range('L89').Value = 95
$ws.Range('L93').Value = 39
$ws.Range('L96').Value = 72
$ws.Range('L99').Value = 117
$ws.Range('L100').Value = 10
$ws.Range('I101').Value = 26303
$ws.Range('L101').Value = 7398

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L6').Value = 25
$ws.Range('L7').Value = 72

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 76
$ws.Range('L6').Value = 72
$ws.Range('L7').Value = 265

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 73

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('L2').Value = 13
$ws.Range('L7').Value = 42

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L6').Value = 79
$ws.Range('L7').Value = 152

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L4').Value = 16
$ws.Range('L7').Value = 376

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L6').Value = 39
$ws.Range('L7').Value = 101

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L3').Value = 14
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L2').Value = 16
$ws.Range('L7').Value = 39

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 60
$ws.Range('L7').Value = 232

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('L5').Value = 4
$ws.Range('L6').Value = 9

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L2').Value = 23
$ws.Range('L7').Value = 48

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L2').Value = 33
$ws.Range('L3').Value = 25
$ws.Range('L6').Value = 29
$ws.Range('L7').Value = 101

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L6').Value = 13
$ws.Range('L7').Value = 66

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L2').Value = 23
$ws.Range('L7').Value = 75

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 72

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 73
$ws.Range('L6').Value = 43
$ws.Range('L7').Value = 204

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L4').Value = 15
$ws.Range('L6').Value = 53
$ws.Range('L7').Value = 188

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L3').Value = 20
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('L6').Value = 12
$ws.Range('L7').Value = 39

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('L3').Value = 4
$ws.Range('L7').Value = 10

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 70
$ws.Range('L3').Value = 77
$ws.Range('L7').Value = 239

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('L3').Value = 22
$ws.Range('L6').Value = 7
$ws.Range('L7').Value = 41

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 51

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L3').Value = 39
$ws.Range('L7').Value = 131

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L2').Value = 23
$ws.Range('L7').Value = 58

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L3').Value = 23
$ws.Range('L4').Value = 16
$ws.Range('L7').Value = 95

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L6').Value = 18
$ws.Range('L7').Value = 76

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 34
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L5').Value = 2
$ws.Range('L7').Value = 84

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('L3').Value = 8
$ws.Range('L7').Value = 32

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L2').Value = 12
$ws.Range('L7').Value = 59

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 119
$ws.Range('L3').Value = 161
$ws.Range('L5').Value = 8
$ws.Range('L7').Value = 386

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range('L6').Value = 2
$ws.Range('L7').Value = 6

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('L2').Value = 4
$ws.Range('L7').Value = 17
